# Auto-generated edit script: updates price/profit columns (H:N)
# for specific Leve rows across multiple worksheets, per the
# scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 782.6429000000001
$ws.Range("I107").Value = 912.9
$ws.Range("J107").Value = 457
$ws.Range("K107").Value = 912.9
$ws.Range("L107").Value = 457
$ws.Range("M107").Value = 1007.1
$ws.Range("N107").Value = -4297
# Row 132
$ws.Range("H132").Value = 8339801.5
$ws.Range("I132").Value = 10006692
$ws.Range("K132").Value = 30020076
$ws.Range("M132").Value = -30017546
# Row 137
$ws.Range("H137").Value = 6672606
$ws.Range("I137").Value = 11118865
$ws.Range("J137").Value = 3216.6667
$ws.Range("K137").Value = 33356595
$ws.Range("L137").Value = 9650.000100000001
$ws.Range("M137").Value = -33354045
$ws.Range("N137").Value = -14750.0001
# Row 138
$ws.Range("H138").Value = 2365.7378
$ws.Range("I138").Value = 1156.1143
$ws.Range("J138").Value = 3994.077
$ws.Range("K138").Value = 3468.3429
$ws.Range("L138").Value = 11982.231
$ws.Range("M138").Value = 1671.6571
$ws.Range("N138").Value = -22262.231

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 10001279
$ws.Range("I2").Value = 11364536
$ws.Range("J2").Value = 4060
$ws.Range("K2").Value = 11364536
$ws.Range("L2").Value = 4060
$ws.Range("M2").Value = -11364423
$ws.Range("N2").Value = -4286
# Row 32
$ws.Range("H32").Value = 2335.81
$ws.Range("I32").Value = 1601.216
$ws.Range("J32").Value = 7722.8335
$ws.Range("K32").Value = 1601.216
$ws.Range("L32").Value = 7722.8335
$ws.Range("M32").Value = -1314.216
$ws.Range("N32").Value = -8296.833500000001
# Row 61
$ws.Range("H61").Value = 2776
$ws.Range("I61").Value = 1535.4667
$ws.Range("J61").Value = 3662.0952
$ws.Range("K61").Value = 1535.4667
$ws.Range("L61").Value = 3662.0952
$ws.Range("M61").Value = -1323.4667
$ws.Range("N61").Value = -4086.0952
# Row 116
$ws.Range("H116").Value = 10001279
$ws.Range("I116").Value = 11364536
$ws.Range("J116").Value = 4060
$ws.Range("K116").Value = 11364536
$ws.Range("L116").Value = 4060
$ws.Range("M116").Value = -11362242
$ws.Range("N116").Value = -8648
# Row 122
$ws.Range("H122").Value = 2598.2188
$ws.Range("I122").Value = 1884.9048
$ws.Range("K122").Value = 5654.7144
$ws.Range("M122").Value = -3204.7144
# Row 132
$ws.Range("H132").Value = 27780544
$ws.Range("I132").Value = 41668708
$ws.Range("J132").Value = 4216.5
$ws.Range("K132").Value = 125006124
$ws.Range("L132").Value = 12649.5
$ws.Range("M132").Value = -125003594
$ws.Range("N132").Value = -17709.5
# Row 136
$ws.Range("H136").Value = 2776
$ws.Range("I136").Value = 1535.4667
$ws.Range("J136").Value = 3662.0952
$ws.Range("K136").Value = 4606.4001
$ws.Range("L136").Value = 10986.2856
$ws.Range("M136").Value = -2056.4001
$ws.Range("N136").Value = -16086.2856

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 10001279
$ws.Range("I3").Value = 11364536
$ws.Range("J3").Value = 4060
$ws.Range("K3").Value = 11364536
$ws.Range("L3").Value = 4060
$ws.Range("M3").Value = -11364422
$ws.Range("N3").Value = -4288
# Row 86
$ws.Range("H86").Value = 35898.668
$ws.Range("I86").Value = 1208
$ws.Range("J86").Value = 105280
$ws.Range("K86").Value = 1208
$ws.Range("L86").Value = 105280
$ws.Range("M86").Value = -85
$ws.Range("N86").Value = -107526
# Row 89
$ws.Range("H89").Value = 35898.668
$ws.Range("I89").Value = 1208
$ws.Range("J89").Value = 105280
$ws.Range("K89").Value = 6040
$ws.Range("L89").Value = 526400
$ws.Range("M89").Value = -424
$ws.Range("N89").Value = -537632
# Row 99
$ws.Range("H99").Value = 2519.9285
$ws.Range("I99").Value = 1727.9
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 1727.9
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -229.9000000000001
$ws.Range("N99").Value = -7496
# Row 105
$ws.Range("H105").Value = 1770.92
$ws.Range("I105").Value = 1551.3334
$ws.Range("J105").Value = 2100.3
$ws.Range("K105").Value = 1551.3334
$ws.Range("L105").Value = 2100.3
$ws.Range("M105").Value = 195.6666
$ws.Range("N105").Value = -5594.3

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2275038.8
$ws.Range("I31").Value = 2327714
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 2327714
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -2327419
$ws.Range("N31").Value = -10590
# Row 34
$ws.Range("H34").Value = 2275038.8
$ws.Range("I34").Value = 2327714
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 2327714
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -2327512
$ws.Range("N34").Value = -10404
# Row 132
$ws.Range("H132").Value = 2851.7666
$ws.Range("I132").Value = 2033.65
$ws.Range("J132").Value = 4488
$ws.Range("K132").Value = 6100.950000000001
$ws.Range("L132").Value = 13464
$ws.Range("M132").Value = -3570.950000000001
$ws.Range("N132").Value = -18524
# Row 134
$ws.Range("H134").Value = 1389.5405
$ws.Range("I134").Value = 876.3043
$ws.Range("J134").Value = 2232.7144
$ws.Range("K134").Value = 2628.9129
$ws.Range("L134").Value = 6698.1432
$ws.Range("M134").Value = -93.91290000000026
$ws.Range("N134").Value = -11768.1432

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 101
$ws.Range("H101").Value = 7257.25
$ws.Range("J101").Value = 7257.25
$ws.Range("L101").Value = 21771.75
$ws.Range("N101").Value = -26639.75
# Row 117
$ws.Range("H117").Value = 1265.5714
$ws.Range("J117").Value = 1376.5
$ws.Range("L117").Value = 4129.5
$ws.Range("N117").Value = -11013.5
# Row 118
$ws.Range("H118").Value = 3626.4119
$ws.Range("I118").Value = 1009.3333
$ws.Range("J118").Value = 3879.6775
$ws.Range("K118").Value = 3027.9999
$ws.Range("L118").Value = 11639.0325
$ws.Range("M118").Value = -1784.9999
$ws.Range("N118").Value = -14125.0325

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 19999.5
$ws.Range("I15").Value = 20000
$ws.Range("K15").Value = 20000
$ws.Range("M15").Value = -19712
# Row 57
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 6000
$ws.Range("J57").Value = 14000
$ws.Range("K57").Value = 6000
$ws.Range("L57").Value = 14000
$ws.Range("M57").Value = -5180
$ws.Range("N57").Value = -15640
# Row 70
$ws.Range("H70").Value = 4538.8184
$ws.Range("I70").Value = 4569.625
$ws.Range("J70").Value = 4456.6665
$ws.Range("K70").Value = 4569.625
$ws.Range("L70").Value = 4456.6665
$ws.Range("M70").Value = -4299.625
$ws.Range("N70").Value = -4996.6665
# Row 73
$ws.Range("H73").Value = 4538.8184
$ws.Range("I73").Value = 4569.625
$ws.Range("J73").Value = 4456.6665
$ws.Range("K73").Value = 4569.625
$ws.Range("L73").Value = 4456.6665
$ws.Range("M73").Value = -3633.625
$ws.Range("N73").Value = -6328.6665
# Row 81
$ws.Range("H81").Value = 19999.5
$ws.Range("I81").Value = 20000
$ws.Range("K81").Value = 20000
$ws.Range("M81").Value = -19002
# Row 84
$ws.Range("H84").Value = 19999.5
$ws.Range("I84").Value = 20000
$ws.Range("K84").Value = 60000
$ws.Range("M84").Value = -55008
# Row 132
$ws.Range("H132").Value = 2983.6592
$ws.Range("I132").Value = 2513.8572
$ws.Range("J132").Value = 3805.8125
$ws.Range("K132").Value = 7541.571599999999
$ws.Range("L132").Value = 11417.4375
$ws.Range("M132").Value = -5011.571599999999
$ws.Range("N132").Value = -16477.4375

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 17666.666
$ws.Range("J6").Value = 17666.666
$ws.Range("L6").Value = 17666.666
$ws.Range("N6").Value = -17890.666
# Row 16
$ws.Range("H16").Value = 1616.8462
$ws.Range("I16").Value = 891
$ws.Range("J16").Value = 3250
$ws.Range("K16").Value = 891
$ws.Range("L16").Value = 3250
$ws.Range("M16").Value = -721
$ws.Range("N16").Value = -3590

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1536.9375
$ws.Range("I113").Value = 656
$ws.Range("J113").Value = 2222.111
$ws.Range("K113").Value = 1968
$ws.Range("L113").Value = 6666.333
$ws.Range("M113").Value = 202
$ws.Range("N113").Value = -11006.333
# Row 126
$ws.Range("H126").Value = 4002109
$ws.Range("I126").Value = 1679.8948
$ws.Range("J126").Value = 16670134
$ws.Range("K126").Value = 5039.6844
$ws.Range("L126").Value = 50010402
$ws.Range("M126").Value = -2569.6844
$ws.Range("N126").Value = -50015342
# Row 136
$ws.Range("H136").Value = 1386.375
$ws.Range("I136").Value = 738.65
$ws.Range("J136").Value = 4625
$ws.Range("K136").Value = 2215.95
$ws.Range("L136").Value = 13875
$ws.Range("M136").Value = 334.0500000000002
$ws.Range("N136").Value = -18975
